$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.342.12'
$ws.Range("E2").Value = '  +4.14%  '
$ws.Range("D3").Value = '3.637.84'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '202.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +11.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '568.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.50%  '
$ws.Range("D7").Value = '3.625.35'
$ws.Range("E7").Value = '  +3.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.618'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.37%  '
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.678'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.75'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.154'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000291'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +16.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.75%  '
$ws.Range("D15").Value = '4.212.42'
$ws.Range("E15").Value = '  +3.01%  '
$ws.Range("D16").Value = '3.640.24'
$ws.Range("E16").Value = '  +3.42%  '
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").Value = '68.182.91'
$ws.Range("E18").Value = '  +4.18%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.50%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '401.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +27.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.80%  '
$ws.Range("E27").Value = '  +2.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.84'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +21.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '692.53'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +14.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.23'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.60%  '
$ws.Range("E35").Value = '  +4.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '64.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '42.71'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.425'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +15.44%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '0.0₃0779'
$ws.Range("E40").Value = '  +6.08%  '
$ws.Range("E41").Value = '  +7.80%  '
$ws.Range("D42").Value = '3.246.02'
$ws.Range("E42").Value = '  +14.60%  '
$ws.Range("E43").Value = '  +13.97%  '
$ws.Range("E44").Value = '  +17.56%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +38.98%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.997'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0418'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.38%  '
$ws.Range("E50").Value = '  +2.06%  '
$ws.Range("E51").Value = '  +5.48%  '
